$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.680.16"
$ws.Range("E2").Value = "  +0.12%  "

# Row 3
$ws.Range("D3").Value = "2.224.10"
$ws.Range("E3").Value = "  +2.21%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("D5").Value = "'270.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.05%  "

# Row 6
$ws.Range("D6").Value = "'93.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +15.33%  "

# Row 7
$ws.Range("D7").Value = "'0.628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.25%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("D9").Value = "'0.621"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.61%  "

# Row 10
$ws.Range("D10").Value = "'45.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.05%  "

# Row 11
$ws.Range("D11").Value = "'0.0977"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.55%  "

# Row 12
$ws.Range("D12").Value = "'8.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +18.68%  "

# Row 13
$ws.Range("D13").Value = "'0.105"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.66%  "

# Row 14
$ws.Range("D14").Value = "2.558.94"
$ws.Range("E14").Value = "  +1.77%  "

# Row 15
$ws.Range("D15").Value = "'15.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.53%  "

# Row 16
$ws.Range("D16").Value = "'0.806"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.52%  "

# Row 17
$ws.Range("D17").Value = "2.221.11"
$ws.Range("E17").Value = "  +2.05%  "

# Row 18
$ws.Range("D18").Value = "43.645.14"
$ws.Range("E18").Value = "  +0.53%  "

# Row 19
$ws.Range("D19").Value = "'0.0000106"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.49%  "

# Row 20
$ws.Range("D20").Value = "'6.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.56%  "

# Row 21
$ws.Range("D21").Value = "'70.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.23%  "

# Row 22
$ws.Range("D22").Value = "'2.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.24%  "

# Row 23
$ws.Range("D23").Value = "'233.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.72%  "

# Row 24
$ws.Range("D24").Value = "'9.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.50%  "

# Row 25
$ws.Range("E25").Value = "  +0.13%  "

# Row 26
$ws.Range("D26").Value = "'11.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.97%  "

# Row 27
$ws.Range("E27").Value = "  +13.37%  "

# Row 28
$ws.Range("D28").Value = "'42.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.71%  "

# Row 29
$ws.Range("D29").Value = "'3.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.20%  "

# Row 30
$ws.Range("E30").Value = "  +2.09%  "

# Row 31
$ws.Range("D31").Value = "'172.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.18%  "

# Row 32
$ws.Range("D32").Value = "'0.0919"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.56%  "

# Row 33
$ws.Range("D33").Value = "'20.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.89%  "

# Row 34
$ws.Range("D34").Value = "'5.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.28%  "

# Row 35
$ws.Range("D35").Value = "'0.124"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.99%  "

# Row 36
$ws.Range("D36").Value = "'0.113"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.83%  "

# Row 37
$ws.Range("D37").Value = "'0.0351"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.04%  "

# Row 38
$ws.Range("D38").Value = "'4.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.16%  "

# Row 39
$ws.Range("D39").Value = "'3.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +24.58%  "

# Row 40
$ws.Range("D40").Value = "'12.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.20%  "

# Row 41
$ws.Range("D41").Value = "'0.221"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.56%  "

# Row 42
$ws.Range("E42").Value = "  +3.72%  "

# Row 43
$ws.Range("D43").Value = "'63.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.12%  "

# Row 44
$ws.Range("D44").Value = "'5.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.67%  "

# Row 45
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'8.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.87%  "

# Row 46
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.0987"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.99%  "

# Row 47
$ws.Range("D47").Value = "'99.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.94%  "

# Row 48
$ws.Range("E48").Value = "  +4.49%  "

# Row 49
$ws.Range("D49").Value = "'1.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.61%  "

# Row 50
$ws.Range("D50").Value = "'0.445"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.53%  "

# Row 51
$ws.Range("E51").Value = "  -5.47%  "
